# SMKCO7004.xlsx update -> 2.22.0-SNAPSHOT
#
# The "Danish term" column is dropped, the trailing "(system)API Search
# Form" / "(system)API Search Form Sorting" columns are replaced by a
# single "(system)API Property Mapping" column (placed before the
# "(system)API Search Criteria Mapping" column), row 1 is given a fixed
# 14.25pt height, and the header row is selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift B1:P1 left by one (dropping "Danish term" from B1) and rebuild the
# last two header cells so that O1/P1 hold the new layout; then drop the
# now-unused Q1:R1 cells entirely.
$ws.Range("B1").Value = "Imported from CMS (Y/N)"
$ws.Range("C1").Value = "Optionality  (O/M)CO is imported from CMS"
$ws.Range("D1").Value = "Optionality  (O/M)CO is created in CS"
$ws.Range("E1").Value = "Editing  (E/R/H)CO is imported from CMS"
$ws.Range("F1").Value = "Editing  (E/R/H)CO is created in CS"
$ws.Range("G1").Value = "Data Type"
$ws.Range("H1").Value = "Controlled  vocabulary"
$ws.Range("I1").Value = "Searchable"
$ws.Range("J1").Value = "Description"
$ws.Range("K1").Value = "Sample Data"
$ws.Range("L1").Value = "Comment"
$ws.Range("M1").Value = "(system)CS URI"
$ws.Range("N1").Value = "(system)CS Definition Property Mapping"
$ws.Range("O1").Value = "(system)API Property Mapping"
$ws.Range("P1").Value = "(system)API Search Criteria Mapping"
$ws.Range("Q1:R1").ClearContents()

# Row 1 now renders slightly shorter (custom height) in the refreshed export.
$ws.Rows.Item(1).RowHeight = 14.25

# The refreshed export leaves the header row selected instead of E9.
$ws.Rows.Item(1).Select() | Out-Null
